# Auto-generated edit script: updates cached market-data values in the
# Leve profit sheets to match the scheduled runner's refreshed snapshot.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2456.2
$ws.Range("I51").Value = 2447.4211
$ws.Range("J51").Value = 2471.3635
$ws.Range("K51").Value = 2447.4211
$ws.Range("L51").Value = 2471.3635
$ws.Range("M51").Value = -1963.4211
$ws.Range("N51").Value = -3439.3635
$ws.Range("H53").Value = 171.44444
$ws.Range("I53").Value = 127.38461
$ws.Range("J53").Value = 286
$ws.Range("K53").Value = 127.38461
$ws.Range("L53").Value = 286
$ws.Range("M53").Value = 509.61539
$ws.Range("N53").Value = -1560
$ws.Range("H70").Value = 50094.617
$ws.Range("I70").Value = 1243.1428
$ws.Range("J70").Value = 74520.36
$ws.Range("K70").Value = 3729.4284
$ws.Range("L70").Value = 223561.08
$ws.Range("M70").Value = -3459.4284
$ws.Range("N70").Value = -224101.08
$ws.Range("H73").Value = 50094.617
$ws.Range("I73").Value = 1243.1428
$ws.Range("J73").Value = 74520.36
$ws.Range("K73").Value = 3729.4284
$ws.Range("L73").Value = 223561.08
$ws.Range("M73").Value = -2793.4284
$ws.Range("N73").Value = -225433.08
$ws.Range("H138").Value = 6026678
$ws.Range("I138").Value = 1287.4642
$ws.Range("J138").Value = 9094150
$ws.Range("K138").Value = 3862.3926
$ws.Range("L138").Value = 27282450
$ws.Range("M138").Value = 1277.6074
$ws.Range("N138").Value = -27292730

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 3667.0754
$ws.Range("I61").Value = 2105.1667
$ws.Range("K61").Value = 2105.1667
$ws.Range("M61").Value = -1893.1667
$ws.Range("H62").Value = 80000
$ws.Range("J62").Value = 80000
$ws.Range("L62").Value = 80000
$ws.Range("N62").Value = -81248
$ws.Range("H65").Value = 80000
$ws.Range("J65").Value = 80000
$ws.Range("L65").Value = 240000
$ws.Range("N65").Value = -246240
$ws.Range("H74").Value = 13214.6875
$ws.Range("J74").Value = 31563.834
$ws.Range("L74").Value = 31563.834
$ws.Range("N74").Value = -33311.834
$ws.Range("H77").Value = 13214.6875
$ws.Range("J77").Value = 31563.834
$ws.Range("L77").Value = 157819.17
$ws.Range("N77").Value = -166555.17
$ws.Range("H97").Value = 1739.3334
$ws.Range("I97").Value = 1564.5625
$ws.Range("J97").Value = 2088.875
$ws.Range("K97").Value = 1564.5625
$ws.Range("L97").Value = 2088.875
$ws.Range("M97").Value = -1068.5625
$ws.Range("N97").Value = -3080.875
$ws.Range("H122").Value = 1421.9
$ws.Range("I122").Value = 1164
$ws.Range("J122").Value = 2453.5
$ws.Range("K122").Value = 3492
$ws.Range("L122").Value = 7360.5
$ws.Range("M122").Value = -1042
$ws.Range("N122").Value = -12260.5
$ws.Range("H132").Value = 2979.6562
$ws.Range("I132").Value = 2852.7693
$ws.Range("K132").Value = 8558.3079
$ws.Range("M132").Value = -6028.3079
$ws.Range("H136").Value = 3667.0754
$ws.Range("I136").Value = 2105.1667
$ws.Range("K136").Value = 6315.500100000001
$ws.Range("M136").Value = -3765.500100000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1673.2142
$ws.Range("I86").Value = 1640.4
$ws.Range("J86").Value = 1755.25
$ws.Range("K86").Value = 1640.4
$ws.Range("L86").Value = 1755.25
$ws.Range("M86").Value = -517.4000000000001
$ws.Range("N86").Value = -4001.25
$ws.Range("H89").Value = 1673.2142
$ws.Range("I89").Value = 1640.4
$ws.Range("J89").Value = 1755.25
$ws.Range("K89").Value = 8202
$ws.Range("L89").Value = 8776.25
$ws.Range("M89").Value = -2586
$ws.Range("N89").Value = -20008.25
$ws.Range("H107").Value = 1094.1904
$ws.Range("I107").Value = 1062.421
$ws.Range("K107").Value = 1062.421
$ws.Range("M107").Value = 857.579

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2540.4167
$ws.Range("I99").Value = 2444.5
$ws.Range("J99").Value = 3020
$ws.Range("K99").Value = 2444.5
$ws.Range("L99").Value = 3020
$ws.Range("M99").Value = -946.5
$ws.Range("N99").Value = -6016
$ws.Range("H126").Value = 2540.4167
$ws.Range("I126").Value = 2444.5
$ws.Range("J126").Value = 3020
$ws.Range("K126").Value = 7333.5
$ws.Range("L126").Value = 9060
$ws.Range("M126").Value = -4863.5
$ws.Range("N126").Value = -14000
$ws.Range("H141").Value = 255754.9
$ws.Range("J141").Value = 329788.12
$ws.Range("L141").Value = 329788.12
$ws.Range("N141").Value = -340148.12

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 629.5454999999999
$ws.Range("J5").Value = 1100
$ws.Range("L5").Value = 3300
$ws.Range("N5").Value = -3524
$ws.Range("H135").Value = 629.5454999999999
$ws.Range("J135").Value = 1100
$ws.Range("L135").Value = 9900
$ws.Range("N135").Value = -14970

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 126.958336
$ws.Range("I2").Value = 91.31579000000001
$ws.Range("K2").Value = 91.31579000000001
$ws.Range("M2").Value = 21.68420999999999
$ws.Range("H21").Value = 9999.5
$ws.Range("J21").Value = 9999.5
$ws.Range("L21").Value = 9999.5
$ws.Range("N21").Value = -10345.5
$ws.Range("H30").Value = 9999.5
$ws.Range("J30").Value = 9999.5
$ws.Range("L30").Value = 9999.5
$ws.Range("N30").Value = -10209.5
$ws.Range("H122").Value = 2960.8823
$ws.Range("I122").Value = 2550.9285
$ws.Range("J122").Value = 4874
$ws.Range("K122").Value = 7652.7855
$ws.Range("L122").Value = 14622
$ws.Range("M122").Value = -5202.7855
$ws.Range("N122").Value = -19522
$ws.Range("H132").Value = 3354.647
$ws.Range("I132").Value = 2746.2856
$ws.Range("J132").Value = 6193.6665
$ws.Range("K132").Value = 8238.856800000001
$ws.Range("L132").Value = 18580.9995
$ws.Range("M132").Value = -5708.856800000001
$ws.Range("N132").Value = -23640.9995

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 1084
$ws.Range("J22").Value = 3543.5557
$ws.Range("K22").Value = 1084
$ws.Range("L22").Value = 3543.5557
$ws.Range("M22").Value = -789
$ws.Range("N22").Value = -4133.5557
$ws.Range("I27").Value = 1084
$ws.Range("J27").Value = 3543.5557
$ws.Range("K27").Value = 1084
$ws.Range("L27").Value = 3543.5557
$ws.Range("M27").Value = -977
$ws.Range("N27").Value = -3757.5557
$ws.Range("H48").Value = 23999.5
$ws.Range("I48").Value = 23999.5
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 23999.5
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -23338.5
$ws.Range("N48").ClearContents()
$ws.Range("H61").Value = 3903.8147
$ws.Range("I61").Value = 3654.7144
$ws.Range("K61").Value = 3654.7144
$ws.Range("M61").Value = -3452.7144
$ws.Range("H74").Value = 29999
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 29999
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H93").Value = 1435.3478
$ws.Range("I93").Value = 1323
$ws.Range("J93").Value = 1839.8
$ws.Range("K93").Value = 1323
$ws.Range("L93").Value = 1839.8
$ws.Range("M93").Value = -75
$ws.Range("N93").Value = -4335.8
$ws.Range("H100").Value = 9825.125
$ws.Range("I100").Value = 4086
$ws.Range("K100").Value = 4086
$ws.Range("M100").Value = -3545
$ws.Range("H113").Value = 3903.8147
$ws.Range("I113").Value = 3654.7144
$ws.Range("K113").Value = 3654.7144
$ws.Range("M113").Value = -1484.7144
$ws.Range("H122").Value = 5599.467
$ws.Range("I122").Value = 5374.25
$ws.Range("J122").Value = 5856.857
$ws.Range("K122").Value = 16122.75
$ws.Range("L122").Value = 17570.571
$ws.Range("M122").Value = -13672.75
$ws.Range("N122").Value = -22470.571

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 644.4643
$ws.Range("I100").Value = 633
$ws.Range("J100").Value = 673.125
$ws.Range("K100").Value = 1266
$ws.Range("L100").Value = 1346.25
$ws.Range("M100").Value = -725
$ws.Range("N100").Value = -2428.25
$ws.Range("H107").Value = 53369.105
$ws.Range("I107").Value = 730.2308
$ws.Range("J107").Value = 167420
$ws.Range("K107").Value = 2190.6924
$ws.Range("L107").Value = 502260
$ws.Range("M107").Value = -270.6923999999999
$ws.Range("N107").Value = -506100
$ws.Range("H113").Value = 1773.1852
$ws.Range("I113").Value = 1115.08
$ws.Range("K113").Value = 3345.24
$ws.Range("M113").Value = -1175.24
$ws.Range("H132").Value = 1774.5454
$ws.Range("I132").Value = 1676.2069
$ws.Range("K132").Value = 5028.620699999999
$ws.Range("M132").Value = -2498.620699999999

